# aggiornamento fino a 27/05
# Appends new daily rows (256-269, dates 2021-05-14 .. 2021-05-27) to the
# existing Castelvetro report sheet, extending the data range from
# A1:D255 to A1:D269.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 255
$firstNewRow = 256
$lastNewRow = 269

# Propagate the formatting (date number format / border / alignment on
# column A, plain numeric formatting on B:D) from the last existing row
# down across all of the freshly appended rows before writing values.
$ws.Range("A" + $lastRow + ":D" + $lastRow).Copy($ws.Range("A" + $firstNewRow + ":D" + $lastNewRow))

$data = @(
  @(256,44330,0,12,106.2981663566303),
  @(257,44331,1,10,88.58180529719195),
  @(258,44332,2,8,70.86544423775356),
  @(259,44333,0,8,70.86544423775356),
  @(260,44334,1,7,62.00726370803437),
  @(261,44335,1,8,70.86544423775356),
  @(262,44336,0,5,44.29090264859597),
  @(263,44337,3,8,70.86544423775356),
  @(264,44338,0,7,62.00726370803437),
  @(265,44339,0,5,44.29090264859597),
  @(266,44340,0,5,44.29090264859597),
  @(267,44341,0,4,35.43272211887678),
  @(268,44342,1,4,35.43272211887678),
  @(269,44343,2,6,53.14908317831517)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
